# "Added last minute updates"
#
# 1) Give the first paragraph (the hidden **ID__...__ID** marker line) a
#    paragraph border whose only property is 5pt spacing on all four sides,
#    and bump its left indent from 120 twips (6pt) to 225 twips (11.25pt).
# 2) Fix the placeholder id text itself (5304_topic_4 -> SUBPART_5304_4) and
#    drop the trailing run that held nothing but a single space.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# --- paragraph formatting -------------------------------------------------
# LeftIndent is expressed in points on the Word object model (1pt = 20 twips),
# so 225 twips == 11.25pt.
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Setting only the DistanceFrom* properties (without flipping Borders.Enable)
# yields a <w:pBdr> whose edges carry just w:space, matching the diff.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# --- text content ----------------------------------------------------------
# Update the placeholder id in place (formatting of the run is preserved).
$findRange = $p1.Range.Duplicate
[void]$findRange.Find.Execute("**ID__AFFARS_5304_topic_4__ID**", $false, $false, `
    $false, $false, $false, $true, 1, $false, `
    "**ID__AFFARS_SUBPART_5304_4__ID**", 2)

# The paragraph's second run was nothing but a trailing space; remove it
# (found relative to the paragraph end so it is robust to the text edit
# above shifting absolute character offsets).
$pEnd = $p1.Range.End
$trailing = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailing.Text -eq " ") {
    $trailing.Delete()
}
